$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Columns("B").Delete(-4159)
